$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Fetching Movie Data from api",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Fetching Movie Data from api. Movies data are sorted according to released date.",
    2
)
